$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 83

# Text columns: force text so Excel doesn't auto-convert "2025-02-22" to a
# date serial or "07" to the number 7, then clear the resulting number
# format override so the cell keeps the sheet's default (unstyled) look.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-02-22"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "13:35:43"
$ws.Cells.Item($row, 2).ClearFormats()

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "Saturday"
$ws.Cells.Item($row, 3).ClearFormats()

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "07"
$ws.Cells.Item($row, 4).ClearFormats()

# Numeric columns
$ws.Cells.Item($row, 5).Value = 130532
$ws.Cells.Item($row, 6).Value = 141569
$ws.Cells.Item($row, 7).Value = 172290
$ws.Cells.Item($row, 8).Value = 157897
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 146509
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 193511
$ws.Cells.Item($row, 14).Value = 115395
$ws.Cells.Item($row, 15).Value = 46264
$ws.Cells.Item($row, 16).Value = 29305
$ws.Cells.Item($row, 17).Value = 68153
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 48445
$ws.Cells.Item($row, 20).Value = -1
